$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for data rows 2-11 from 45175 to 45183
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 45183
}
